$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,10
$row2[0,0] = -17.68044194203702
$row2[0,1] = 0.3409555837965801
$row2[0,2] = -17.68044194203702
$row2[0,3] = -17.68044194203702
$row2[0,4] = -17.68044194203702
$row2[0,5] = -17.68044194203702
$row2[0,6] = -17.68044194203702
$row2[0,7] = -17.68044194203702
$row2[0,8] = -17.68044194203702
$row2[0,9] = -17.68044194203702
$ws.Range("B2:K2").Value = $row2

$row3 = New-Object "object[,]" 1,10
$row3[0,0] = -17.68044194203702
$row3[0,1] = -17.68044194203702
$row3[0,2] = -17.68044194203702
$row3[0,3] = -17.68044194203702
$row3[0,4] = -17.68044194203702
$row3[0,5] = -17.68044194203702
$row3[0,6] = -17.68044194203702
$row3[0,7] = 0.8931201923372944
$row3[0,8] = -17.68044194203702
$row3[0,9] = -17.68044194203702
$ws.Range("B3:K3").Value = $row3

$row4 = New-Object "object[,]" 1,10
$row4[0,0] = -17.68044194203702
$row4[0,1] = 0.489331400954583
$row4[0,2] = 1.070961548903369
$row4[0,3] = -17.68044194203702
$row4[0,4] = -17.68044194203702
$row4[0,5] = -17.68044194203702
$row4[0,6] = 1.769724206169444
$row4[0,7] = -17.68044194203702
$row4[0,8] = 2.48024813654327
$row4[0,9] = -17.68044194203702
$ws.Range("B4:K4").Value = $row4

$row5 = New-Object "object[,]" 1,10
$row5[0,0] = -17.68044194203702
$row5[0,1] = 0.885853704719568
$row5[0,2] = -17.68044194203702
$row5[0,3] = -17.68044194203702
$row5[0,4] = -17.68044194203702
$row5[0,5] = 3.309388392952367
$row5[0,6] = -17.68044194203702
$row5[0,7] = -17.68044194203702
$row5[0,8] = -17.68044194203702
$row5[0,9] = -17.68044194203702
$ws.Range("B5:K5").Value = $row5

$row6 = New-Object "object[,]" 1,10
$row6[0,0] = -17.68044194203702
$row6[0,1] = -17.68044194203702
$row6[0,2] = -17.68044194203702
$row6[0,3] = -17.68044194203702
$row6[0,4] = -17.68044194203702
$row6[0,5] = -17.68044194203702
$row6[0,6] = -17.68044194203702
$row6[0,7] = -17.68044194203702
$row6[0,8] = -17.68044194203702
$row6[0,9] = -17.68044194203702
$ws.Range("B6:K6").Value = $row6

$row7 = New-Object "object[,]" 1,10
$row7[0,0] = 3.132055813251355
$row7[0,1] = -17.68044194203702
$row7[0,2] = -17.68044194203702
$row7[0,3] = -17.68044194203702
$row7[0,4] = -17.68044194203702
$row7[0,5] = -17.68044194203702
$row7[0,6] = -17.68044194203702
$row7[0,7] = -17.68044194203702
$row7[0,8] = -17.68044194203702
$row7[0,9] = -17.68044194203702
$ws.Range("B7:K7").Value = $row7

$row8 = New-Object "object[,]" 1,10
$row8[0,0] = -17.68044194203702
$row8[0,1] = -17.68044194203702
$row8[0,2] = -17.68044194203702
$row8[0,3] = 1.792266743593451
$row8[0,4] = -17.68044194203702
$row8[0,5] = -17.68044194203702
$row8[0,6] = -17.68044194203702
$row8[0,7] = -17.68044194203702
$row8[0,8] = -17.68044194203702
$row8[0,9] = -17.68044194203702
$ws.Range("B8:K8").Value = $row8

$row9 = New-Object "object[,]" 1,10
$row9[0,0] = 3.489681781226008
$row9[0,1] = -17.68044194203702
$row9[0,2] = -17.68044194203702
$row9[0,3] = -17.68044194203702
$row9[0,4] = -17.68044194203702
$row9[0,5] = -17.68044194203702
$row9[0,6] = -17.68044194203702
$row9[0,7] = -17.68044194203702
$row9[0,8] = -17.68044194203702
$row9[0,9] = -17.68044194203702
$ws.Range("B9:K9").Value = $row9

$row10 = New-Object "object[,]" 1,10
$row10[0,0] = -17.68044194203702
$row10[0,1] = -17.68044194203702
$row10[0,2] = -17.68044194203702
$row10[0,3] = -17.68044194203702
$row10[0,4] = -17.68044194203702
$row10[0,5] = -17.68044194203702
$row10[0,6] = -17.68044194203702
$row10[0,7] = 0.5497896345352274
$row10[0,8] = -17.68044194203702
$row10[0,9] = 2.058812918946868
$ws.Range("B10:K10").Value = $row10

$row11 = New-Object "object[,]" 1,10
$row11[0,0] = -17.68044194203702
$row11[0,1] = -17.68044194203702
$row11[0,2] = -17.68044194203702
$row11[0,3] = 1.811021702435366
$row11[0,4] = -17.68044194203702
$row11[0,5] = 1.861904160160017
$row11[0,6] = -17.68044194203702
$row11[0,7] = -17.68044194203702
$row11[0,8] = -17.68044194203702
$row11[0,9] = 1.417664696614542
$ws.Range("B11:K11").Value = $row11

$row12 = New-Object "object[,]" 1,10
$row12[0,0] = -17.68044194203702
$row12[0,1] = -17.68044194203702
$row12[0,2] = -17.68044194203702
$row12[0,3] = -17.68044194203702
$row12[0,4] = -17.68044194203702
$row12[0,5] = -17.68044194203702
$row12[0,6] = -17.68044194203702
$row12[0,7] = -17.68044194203702
$row12[0,8] = -17.68044194203702
$row12[0,9] = -17.68044194203702
$ws.Range("B12:K12").Value = $row12

$row13 = New-Object "object[,]" 1,10
$row13[0,0] = -17.68044194203702
$row13[0,1] = -17.68044194203702
$row13[0,2] = -17.68044194203702
$row13[0,3] = 1.470170586472832
$row13[0,4] = -17.68044194203702
$row13[0,5] = -17.68044194203702
$row13[0,6] = -17.68044194203702
$row13[0,7] = -17.68044194203702
$row13[0,8] = 1.375880889995759
$row13[0,9] = 2.405374105952688
$ws.Range("B13:K13").Value = $row13

$row14 = New-Object "object[,]" 1,10
$row14[0,0] = -17.68044194203702
$row14[0,1] = -17.68044194203702
$row14[0,2] = 1.346311117797577
$row14[0,3] = -17.68044194203702
$row14[0,4] = -17.68044194203702
$row14[0,5] = -17.68044194203702
$row14[0,6] = -17.68044194203702
$row14[0,7] = -17.68044194203702
$row14[0,8] = -17.68044194203702
$row14[0,9] = 1.93334212785543
$ws.Range("B14:K14").Value = $row14

$row15 = New-Object "object[,]" 1,10
$row15[0,0] = -17.68044194203702
$row15[0,1] = -17.68044194203702
$row15[0,2] = 0.06526520353412144
$row15[0,3] = -17.68044194203702
$row15[0,4] = -17.68044194203702
$row15[0,5] = -17.68044194203702
$row15[0,6] = -17.68044194203702
$row15[0,7] = -17.68044194203702
$row15[0,8] = -17.68044194203702
$row15[0,9] = -17.68044194203702
$ws.Range("B15:K15").Value = $row15

$row16 = New-Object "object[,]" 1,10
$row16[0,0] = -17.68044194203702
$row16[0,1] = -17.68044194203702
$row16[0,2] = -17.68044194203702
$row16[0,3] = -17.68044194203702
$row16[0,4] = -17.68044194203702
$row16[0,5] = -17.68044194203702
$row16[0,6] = -17.68044194203702
$row16[0,7] = -17.68044194203702
$row16[0,8] = 2.046795508363998
$row16[0,9] = -17.68044194203702
$ws.Range("B16:K16").Value = $row16

$row17 = New-Object "object[,]" 1,10
$row17[0,0] = -17.68044194203702
$row17[0,1] = 0.6838338080497419
$row17[0,2] = 0.4072953978747913
$row17[0,3] = -17.68044194203702
$row17[0,4] = -17.68044194203702
$row17[0,5] = -17.68044194203702
$row17[0,6] = 2.124484704625488
$row17[0,7] = 0.5079971946581766
$row17[0,8] = 2.060427462533169
$row17[0,9] = -17.68044194203702
$ws.Range("B17:K17").Value = $row17

$row18 = New-Object "object[,]" 1,10
$row18[0,0] = -17.68044194203702
$row18[0,1] = -17.68044194203702
$row18[0,2] = -17.68044194203702
$row18[0,3] = -17.68044194203702
$row18[0,4] = -17.68044194203702
$row18[0,5] = -17.68044194203702
$row18[0,6] = 2.032597135536488
$row18[0,7] = 0.1502775035385076
$row18[0,8] = 1.816296780201818
$row18[0,9] = -17.68044194203702
$ws.Range("B18:K18").Value = $row18

$row19 = New-Object "object[,]" 1,10
$row19[0,0] = -17.68044194203702
$row19[0,1] = -17.68044194203702
$row19[0,2] = 2.788354099494711
$row19[0,3] = -17.68044194203702
$row19[0,4] = -17.68044194203702
$row19[0,5] = -17.68044194203702
$row19[0,6] = 1.685369995965336
$row19[0,7] = 1.517573776212793
$row19[0,8] = -17.68044194203702
$row19[0,9] = -17.68044194203702
$ws.Range("B19:K19").Value = $row19

$row20 = New-Object "object[,]" 1,10
$row20[0,0] = -17.68044194203702
$row20[0,1] = 2.809042484999187
$row20[0,2] = 2.603027414598925
$row20[0,3] = -17.68044194203702
$row20[0,4] = 4.321921570259217
$row20[0,5] = -17.68044194203702
$row20[0,6] = 1.320694071939623
$row20[0,7] = 3.496196903133797
$row20[0,8] = -17.68044194203702
$row20[0,9] = 2.016085878153854
$ws.Range("B20:K20").Value = $row20

$row21 = New-Object "object[,]" 1,10
$row21[0,0] = -17.68044194203702
$row21[0,1] = 2.779698946292773
$row21[0,2] = -17.68044194203702
$row21[0,3] = 3.358520571941848
$row21[0,4] = -17.68044194203702
$row21[0,5] = 2.689654762708275
$row21[0,6] = 1.27743378489441
$row21[0,7] = -17.68044194203702
$row21[0,8] = -17.68044194203702
$row21[0,9] = -19.83129332503361
$ws.Range("B21:K21").Value = $row21
